$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 132
$ws.Range("E132").Value = 'Banka 25 TL hesap işletim ücreti alıyor. Ek-1 listesinde hesap işletim ücreti belirtilmemiş ve bu ücretin alınması yasaktır. Tebliğ''de açıkça yasaklanmış bir ücret alınmaktadır.'
$ws.Range("F132").Value = 'Madde 11/Fıkra 1, Madde 9/Fıkra 1, Madde Ek-1/Fıkra Tam Liste'

# Row 133
$ws.Range("E133").Value = 'Banka kendi ATM''sinden para çekme işleminden 5 TL ücret alıyor. Tebliğ''e göre kendi ATM''lerinden para çekme işlemlerinden ücret alınamaz. Bu durum Tebliğ''in açıkça yasakladığı bir durumu içeriyor.'

# Row 134
$ws.Range("D134").Value = 'NA'
$ws.Range("E134").Value = 'Banka maddesi ''Başka Banka ATM Nakit Çekim'' ücretinden bahsediyor. Ek-1 listesinde ''Başka Kuruluş ATM''sinden Yapılan İşlem Ücreti'' yer alıyor. Bu, ücretin izin verildiği anlamına gelir. Ancak, Tebliğ''de bu işlem için bir sayısal limit belirtilmemiştir. Bu nedenle, ihlal olup olmadığını belirlemek mümkün değildir. Bu durumda NA vermek en doğru yaklaşımdır.'
$ws.Range("F134").Value = 'Madde 6/Fıkra 1, Madde 6/Fıkra 2, Madde 9/Fıkra 6, Madde Ek-1/Fıkra Tam Liste'
$ws.Range("G134").Value = 'Gerekli Değil'

# Row 135
$ws.Range("F135").Value = 'Madde 6/Fıkra 1, Madde 6/Fıkra 2, Madde Ek-1/Fıkra Tam Liste'

# Row 136
$ws.Range("E136").Value = 'Banka 15.0 TL şube EFT ücreti alıyor, ancak Tebliğ''e göre 10.0 TL''yi geçemiyor. 5.0 TL fazla ücret alınıyor.'
$ws.Range("F136").Value = 'Madde 11/Fıkra 6, Madde 9/Fıkra 1, Madde Ek-1/Fıkra Tam Liste'

# Row 137
$ws.Range("E137").Value = 'Banka 50.0 TL hesap açılış ücreti alıyor. Tebliğ''in 13. maddesinin 1. fıkrası, mevduat ve katılım fonu hesaplarının açılış işlemlerinde ücret alınamayacağını belirtiyor. Bu durum, Tebliğ''e aykırıdır.'
$ws.Range("F137").Value = 'Madde 12/Fıkra 2, Madde 12/Fıkra 3, Madde 13/Fıkra 1, Madde 6/Fıkra 2, Madde Ek-1/Fıkra Tam Liste'
$ws.Range("G137").Value = 'Hesap açılış işlemlerinde herhangi bir ücret alınmayacaktır.'

# Row 138
$ws.Range("E138").Value = 'Banka 15 TL hesap bakım ücreti alıyor. Tebliğ''e göre mevduat hesaplarından hesap işletim ücreti alınamaz. Bu durum, Tebliğ''in açıkça yasakladığı bir ücretlendirme uygulamasıdır.'
$ws.Range("F138").Value = 'Madde 11/Fıkra 1, Madde 9/Fıkra 1, Madde Ek-1/Fıkra Tam Liste'

# Row 139
$ws.Range("E139").Value = 'Banka SMS bilgilendirme hizmeti için 10 TL ücret alıyor. Tebliğ''de bu hizmetin ücretsiz olması gerekmektedir. Ek-1 listesinde SMS bilgilendirme ücreti yer almamaktadır. Bu nedenle ücret alınması yasaktır.'
$ws.Range("F139").Value = 'Madde 11/Fıkra 6, Madde 9/Fıkra 1, Madde 9/Fıkra 5, Madde Ek-1/Fıkra Tam Liste'

# Row 140
$ws.Range("E140").Value = 'Banka, E-posta bilgilendirme hizmeti için 5.0 TL ücret alıyor. Tebliğ''de bu hizmetin ücretsiz olması gerektiği belirtilmiyor ancak Ek-1 listesinde izin verilen ücret kalemleri arasında yer almıyor. Bu nedenle, ücretsiz olması gereken bir hizmet için ücret alınması Tebliğ''e aykırıdır.'
$ws.Range("F140").Value = 'Madde 11/Fıkra 6, Madde 9/Fıkra 1, Madde 9/Fıkra 5, Madde Ek-1/Fıkra Tam Liste'

# Row 141
$ws.Range("E141").Value = 'Banka, müşteri onayı olmadan SMS bildirimi için ücret alıyor. Tebliğ, müşteri onayı gerektiren işlemler için açıkça izin vermiyor ve bu durum müşteri haklarının ihlali anlamına gelir. Bu, müşteri onayı olmadan ücret alınması yasağına aykırıdır.'
$ws.Range("F141").Value = 'Madde 11/Fıkra 6, Madde 6/Fıkra 2, Madde Ek-1/Fıkra Tam Liste'
$ws.Range("G141").Value = 'Müşteri onayı alınmadan SMS bildirimi gönderilemez. Müşteri onayı alındıktan sonra belirtilen ücret tahsil edilebilir.'

# Row 142
$ws.Range("E142").Value = 'Banka kendi hesabına para yatırma işleminden 3 TL ücret alıyor. Tebliğ''e göre kendi hesabına para yatırma işlemlerinden ücret alınamaz. Bu durum, Tebliğ''in açıkça yasakladığı bir durumu ihlal etmektedir.'
$ws.Range("F142").Value = 'Madde 9/Fıkra 1, Madde Ek-1/Fıkra Tam Liste'
$ws.Range("G142").Value = 'Kendi hesabına para yatırma işlemleri ücretsizdir.'

# Row 143
$ws.Range("E143").Value = 'Banka, sözleşme ilk yıl içinde tekrar basımı için 25 TL ücret alıyor. Tebliğ''e göre sözleşmenin bir örneği ilk yıl ücretsiz verilmesi zorunludur. Bu durum, Tebliğ''in [MADDE 9 - FIKRA 4] hükmünü ihlal etmektedir.'
$ws.Range("F143").Value = 'Madde 9/Fıkra 3, Madde Ek-1/Fıkra Tam Liste'
$ws.Range("G143").Value = 'Sözleşme örneği ilk yıl ücretsizdir.'
